$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- Update the "Date" metadata row (row 8) with the new generation timestamp ---
$ws.Range("B8").Value = "2025-10-02T18:31:12+01:00"

# --- Set the "Case Sensitive" metadata row (row 20) to the literal text "true" ---
# A bare Value = "true" would be auto-typed as a Boolean by Excel, so we enter it
# with a leading apostrophe to force text entry, then re-apply the (unformatted)
# style from a neighboring plain cell so the cell's style index is unchanged.
$srcStyle = $ws.Range("B19")
$dst = $ws.Range("B20")
$dst.Value = "'true"
$srcStyle.Copy()
$dst.PasteSpecial(-4122)  # xlPasteFormats - keep original (unformatted) style, just fix it up after the text entry
$excel.CutCopyMode = $false
